$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.886.77"
$ws.Range("E2").Value = "  +4.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.647.23"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.57"
$ws.Range("E5").Value = "  +6.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.10"
$ws.Range("E6").Value = "  +3.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.609"
$ws.Range("E8").Value = "  +6.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.673.40"
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.87"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  +5.51%  "
$ws.Range("E12").Value = "  +7.14%  "
$ws.Range("E13").Value = "  +3.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.141.61"
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.785.84"
$ws.Range("E15").Value = "  +3.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.97"
$ws.Range("E16").Value = "  +6.09%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("E17").Value = "  +4.84%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.668.61"
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("E19").Value = "  +3.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.21"
$ws.Range("E20").Value = "  +2.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.49"
$ws.Range("E21").Value = "  +3.73%  "
$ws.Range("E22").Value = "  +2.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.82"
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.98"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.441"
$ws.Range("E26").Value = "  +5.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.165"
$ws.Range("E27").Value = "  +1.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.43"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0800"
$ws.Range("E30").Value = "  +9.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.72"
$ws.Range("E32").Value = "  +5.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.29"
$ws.Range("E33").Value = "  +7.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.86"
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.28"
$ws.Range("E35").Value = "  +2.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.10"
$ws.Range("E36").Value = "  +5.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.906"
$ws.Range("E37").Value = "  +7.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.911"
$ws.Range("E38").Value = "  +11.97%  "
$ws.Range("E39").Value = "  +5.66%  "
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("E41").Value = "  +7.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "303.83"
$ws.Range("E42").Value = "  +7.71%  "
$ws.Range("E43").Value = "  +2.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0986"
$ws.Range("E45").Value = "  +4.65%  "
$ws.Range("E46").Value = "  +2.30%  "
$ws.Range("E47").Value = "  +4.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.49"
$ws.Range("E48").Value = "  +3.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "128.57"
$ws.Range("E49").Value = "  +13.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.67"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("E51").Value = "  +5.14%  "
